# Applies the scheduled-runner update to the Sheets workbook (H:N derived-price columns).
$wb = $excel.ActiveWorkbook

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 220041600
$ws.Range("I86").Value = 366733340
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 366733340
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -366732217
$ws.Range("N86").Value = -6246

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 220041600
$ws.Range("I89").Value = 366733340
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 1833666700
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -1833661084
$ws.Range("N89").Value = -31232

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6086.8184
$ws.Range("I132").Value = 5811.3887
$ws.Range("J132").Value = 7326.25
$ws.Range("K132").Value = 17434.1661
$ws.Range("L132").Value = 21978.75
$ws.Range("M132").Value = -14904.1661
$ws.Range("N132").Value = -27038.75

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3185.8333
$ws.Range("I138").Value = 3527.2
$ws.Range("J138").Value = 3096
$ws.Range("K138").Value = 10581.6
$ws.Range("L138").Value = 9288
$ws.Range("M138").Value = -5441.599999999999
$ws.Range("N138").Value = -19568

# ARM!row52
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 94429
$ws.Range("J52").Value = 94429
$ws.Range("L52").Value = 94429
$ws.Range("N52").Value = -95065

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1461.3334
$ws.Range("I122").Value = 1326.2222
$ws.Range("J122").Value = 1866.6666
$ws.Range("K122").Value = 3978.6666
$ws.Range("L122").Value = 5599.9998
$ws.Range("M122").Value = -1528.6666
$ws.Range("N122").Value = -10499.9998

# ARM!row133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2516.6667
$ws.Range("I105").Value = 2120
$ws.Range("K105").Value = 2120
$ws.Range("M105").Value = -373

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1714.75
$ws.Range("I58").Value = 1257.2941
$ws.Range("J58").Value = 2124.0527
$ws.Range("K58").Value = 1257.2941
$ws.Range("L58").Value = 2124.0527
$ws.Range("M58").Value = -1054.2941
$ws.Range("N58").Value = -2530.0527

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 906
$ws.Range("I99").Value = 906
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 906
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 592
$ws.Range("N99").ClearContents()

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1273.1818
$ws.Range("I122").Value = 1188.125
$ws.Range("K122").Value = 3564.375
$ws.Range("M122").Value = -1114.375

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 906
$ws.Range("I126").Value = 906
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2718
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -248
$ws.Range("N126").ClearContents()

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1714.75
$ws.Range("I136").Value = 1257.2941
$ws.Range("J136").Value = 2124.0527
$ws.Range("K136").Value = 3771.8823
$ws.Range("L136").Value = 6372.158100000001
$ws.Range("M136").Value = -1221.8823
$ws.Range("N136").Value = -11472.1581

# CUL!row39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4612.5
$ws.Range("J39").Value = 4612.5
$ws.Range("L39").Value = 13837.5
$ws.Range("N39").Value = -14425.5

# CUL!row60
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2999.0789
$ws.Range("I60").Value = 2995
$ws.Range("J60").Value = 2999.1892
$ws.Range("K60").Value = 8985
$ws.Range("L60").Value = 8997.567599999998
$ws.Range("M60").Value = -8734
$ws.Range("N60").Value = -9499.567599999998

# CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 954.2889
$ws.Range("I68").Value = 687.7778
$ws.Range("J68").Value = 1020.9167
$ws.Range("K68").Value = 2063.3334
$ws.Range("L68").Value = 3062.7501
$ws.Range("M68").Value = -1252.3334
$ws.Range("N68").Value = -4684.7501

# CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 954.2889
$ws.Range("I71").Value = 687.7778
$ws.Range("J71").Value = 1020.9167
$ws.Range("K71").Value = 6190.000199999999
$ws.Range("L71").Value = 9188.2503
$ws.Range("M71").Value = -2134.000199999999
$ws.Range("N71").Value = -17300.2503

# CUL!row107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1406
$ws.Range("I107").Value = 256.6154
$ws.Range("K107").Value = 769.8462000000001
$ws.Range("M107").Value = 1150.1538

# CUL!row110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 15506.75
$ws.Range("J110").Value = 19333.334
$ws.Range("L110").Value = 58000.00199999999
$ws.Range("N110").Value = -66180.00199999999

# CUL!row137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 8933.883
$ws.Range("I137").Value = 11098
$ws.Range("K137").Value = 33294
$ws.Range("M137").Value = -28194

# CUL!row141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 16500
$ws.Range("I141").Value = 19000
$ws.Range("J141").Value = 9000
$ws.Range("K141").Value = 57000
$ws.Range("L141").Value = 27000
$ws.Range("M141").Value = -51820
$ws.Range("N141").Value = -37360

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5745.7896
$ws.Range("I70").Value = 5886.3667
$ws.Range("J70").Value = 5218.625
$ws.Range("K70").Value = 5886.3667
$ws.Range("L70").Value = 5218.625
$ws.Range("M70").Value = -5616.3667
$ws.Range("N70").Value = -5758.625

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5745.7896
$ws.Range("I73").Value = 5886.3667
$ws.Range("J73").Value = 5218.625
$ws.Range("K73").Value = 5886.3667
$ws.Range("L73").Value = 5218.625
$ws.Range("M73").Value = -4950.3667
$ws.Range("N73").Value = -7090.625

# GSM!row97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1515.6666
$ws.Range("I97").Value = 1503.75
$ws.Range("J97").Value = 1611
$ws.Range("K97").Value = 1503.75
$ws.Range("L97").Value = 1611
$ws.Range("M97").Value = -1007.75
$ws.Range("N97").Value = -2603

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1715.6364
$ws.Range("I122").Value = 1815.6666
$ws.Range("J122").Value = 1595.6
$ws.Range("K122").Value = 5446.9998
$ws.Range("L122").Value = 4786.799999999999
$ws.Range("M122").Value = -2996.9998
$ws.Range("N122").Value = -9686.799999999999

# LTW!row93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 7694.8823
$ws.Range("I93").Value = 20878.6
$ws.Range("J93").Value = 2201.6667
$ws.Range("K93").Value = 20878.6
$ws.Range("L93").Value = 2201.6667
$ws.Range("M93").Value = -19630.6
$ws.Range("N93").Value = -4697.6667

# LTW!row133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 55326
$ws.Range("J133").Value = 55326
$ws.Range("L133").Value = 55326
$ws.Range("N133").Value = -60386

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1807.2
$ws.Range("I122").Value = 1674.1111
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 5022.3333
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -2572.3333
$ws.Range("N122").Value = -13915

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2473.0425
$ws.Range("I136").Value = 2318.8572
$ws.Range("J136").Value = 2700.2632
$ws.Range("K136").Value = 6956.571599999999
$ws.Range("L136").Value = 8100.7896
$ws.Range("M136").Value = -4406.571599999999
$ws.Range("N136").Value = -13200.7896
